$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3923.077
$ws.Range("I69").Value = 2000
$ws.Range("J69").Value = 5125
$ws.Range("K69").Value = 6000
$ws.Range("L69").Value = 15375
$ws.Range("M69").Value = -5126
$ws.Range("N69").Value = -17123

$ws.Range("H72").Value = 3923.077
$ws.Range("I72").Value = 2000
$ws.Range("J72").Value = 5125
$ws.Range("K72").Value = 18000
$ws.Range("L72").Value = 46125
$ws.Range("M72").Value = -13632
$ws.Range("N72").Value = -54861

$ws.Range("H113").Value = 5002325
$ws.Range("I113").Value = 6251781.5
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 6251781.5
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = -6248527.5
$ws.Range("N113").Value = -11008

$ws.Range("H116").Value = 8860576
$ws.Range("I116").Value = 4466674
$ws.Range("J116").Value = 15251707
$ws.Range("K116").Value = 4466674
$ws.Range("L116").Value = 15251707
$ws.Range("M116").Value = -4463232
$ws.Range("N116").Value = -15258591

$ws.Range("H132").Value = 4258596
$ws.Range("I132").Value = 967267.9
$ws.Range("J132").Value = 18521018
$ws.Range("K132").Value = 2901803.7
$ws.Range("L132").Value = 55563054
$ws.Range("M132").Value = -2899273.7
$ws.Range("N132").Value = -55568114

$ws.Range("H138").Value = 4468.3447
$ws.Range("I138").Value = 3191.6072
$ws.Range("J138").Value = 5074.2544
$ws.Range("K138").Value = 9574.821599999999
$ws.Range("L138").Value = 15222.7632
$ws.Range("M138").Value = -4434.821599999999
$ws.Range("N138").Value = -25502.7632

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 11824.143
$ws.Range("J37").Value = 17379.75
$ws.Range("L37").Value = 17379.75
$ws.Range("N37").Value = -17925.75

$ws.Range("H74").Value = 8336834
$ws.Range("I74").Value = 1734.4
$ws.Range("J74").Value = 22228666
$ws.Range("K74").Value = 1734.4
$ws.Range("L74").Value = 22228666
$ws.Range("M74").Value = -860.4000000000001
$ws.Range("N74").Value = -22230414

$ws.Range("H77").Value = 8336834
$ws.Range("I77").Value = 1734.4
$ws.Range("J77").Value = 22228666
$ws.Range("K77").Value = 8672
$ws.Range("L77").Value = 111143330
$ws.Range("M77").Value = -4304
$ws.Range("N77").Value = -111152066

$ws.Range("H132").Value = 18205074
$ws.Range("I132").Value = 20437824
$ws.Range("J132").Value = 9617577
$ws.Range("K132").Value = 61313472
$ws.Range("L132").Value = 28852731
$ws.Range("M132").Value = -61310942
$ws.Range("N132").Value = -28857791

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 13081.833
$ws.Range("J50").Value = 13081.833
$ws.Range("L50").Value = 13081.833
$ws.Range("N50").Value = -14331.833

$ws.Range("H51").Value = 24548
$ws.Range("J51").Value = 24548
$ws.Range("L51").Value = 24548
$ws.Range("N51").Value = -26020

$ws.Range("H59").Value = 21250
$ws.Range("J59").Value = 21250
$ws.Range("L59").Value = 21250
$ws.Range("N59").Value = -23540

$ws.Range("H60").Value = 13184.714
$ws.Range("J60").Value = 18180
$ws.Range("L60").Value = 18180
$ws.Range("N60").Value = -19202

$ws.Range("H61").Value = 24548
$ws.Range("J61").Value = 24548
$ws.Range("L61").Value = 24548
$ws.Range("N61").Value = -25244

$ws.Range("H68").Value = 18422.857
$ws.Range("J68").Value = 18422.857
$ws.Range("L68").Value = 18422.857
$ws.Range("N68").Value = -19920.857

$ws.Range("H71").Value = 18422.857
$ws.Range("J71").Value = 18422.857
$ws.Range("L71").Value = 55268.571
$ws.Range("N71").Value = -62756.571

$ws.Range("H74").Value = 27000
$ws.Range("J74").Value = 27000
$ws.Range("L74").Value = 27000
$ws.Range("N74").Value = -28748

$ws.Range("H77").Value = 27000
$ws.Range("J77").Value = 27000
$ws.Range("L77").Value = 81000
$ws.Range("N77").Value = -89736

$ws.Range("H132").Value = 2094.95
$ws.Range("I132").Value = 1369.1818
$ws.Range("J132").Value = 2982
$ws.Range("K132").Value = 4107.5454
$ws.Range("L132").Value = 8946
$ws.Range("M132").Value = -1577.5454
$ws.Range("N132").Value = -14006

$ws.Range("H134").Value = 2003959.1
$ws.Range("I134").Value = 3137.8333
$ws.Range("K134").Value = 9413.499899999999
$ws.Range("M134").Value = -6878.499899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3201.0557
$ws.Range("I39").Value = 499.5
$ws.Range("J39").Value = 3538.75
$ws.Range("K39").Value = 1498.5
$ws.Range("L39").Value = 10616.25
$ws.Range("M39").Value = -1204.5
$ws.Range("N39").Value = -11204.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 83201.28999999999
$ws.Range("I122").Value = 141088.38
$ws.Range("J122").Value = 6018.5
$ws.Range("K122").Value = 423265.14
$ws.Range("L122").Value = 18055.5
$ws.Range("M122").Value = -420815.14
$ws.Range("N122").Value = -22955.5

$ws.Range("H126").Value = 7221.4346
$ws.Range("I126").Value = 18415.5
$ws.Range("K126").Value = 55246.5
$ws.Range("M126").Value = -52776.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 12196516
$ws.Range("I22").Value = 617.6
$ws.Range("J22").Value = 16130677
$ws.Range("K22").Value = 617.6
$ws.Range("L22").Value = 16130677
$ws.Range("M22").Value = -322.6
$ws.Range("N22").Value = -16131267

$ws.Range("H27").Value = 12196516
$ws.Range("I27").Value = 617.6
$ws.Range("J27").Value = 16130677
$ws.Range("K27").Value = 617.6
$ws.Range("L27").Value = 16130677
$ws.Range("M27").Value = -510.6
$ws.Range("N27").Value = -16130891

$ws.Range("H122").Value = 18853216
$ws.Range("I122").Value = 11796747
$ws.Range("J122").Value = 100002600
$ws.Range("K122").Value = 35390241
$ws.Range("L122").Value = 300007800
$ws.Range("M122").Value = -35387791
$ws.Range("N122").Value = -300012700

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 4023
$ws.Range("I100").Value = 5848.222
$ws.Range("K100").Value = 11696.444
$ws.Range("M100").Value = -11155.444

$ws.Range("H113").Value = 131.88889
$ws.Range("I113").Value = 123.375
$ws.Range("J113").Value = 200
$ws.Range("K113").Value = 370.125
$ws.Range("L113").Value = 600
$ws.Range("M113").Value = 1799.875
$ws.Range("N113").Value = -4940

$ws.Range("H132").Value = 1016463.2
$ws.Range("I132").Value = 4022.7036
$ws.Range("J132").Value = 2383257.8
$ws.Range("K132").Value = 12068.1108
$ws.Range("L132").Value = 7149773.399999999
$ws.Range("M132").Value = -9538.110799999999
$ws.Range("N132").Value = -7154833.399999999
